# Insert a new data row before row 436 (shifts existing rows 436..526 down
# to 437..527) and populate it with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 436..526 down to 437..527.
$ws.Rows(436).Insert()

# Populate the newly-inserted row 436 with the new record.
$ws.Cells.Item(436, 1).Value  = 7
$ws.Cells.Item(436, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(436, 3).Value  = "Ñuble"
$ws.Cells.Item(436, 4).Value  = 44641
$ws.Cells.Item(436, 5).Value  = 16
$ws.Cells.Item(436, 6).Value  = 100112004
$ws.Cells.Item(436, 7).Value  = "Cebolla"
$ws.Cells.Item(436, 8).Value  = "Sin especificar"
$ws.Cells.Item(436, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(436, 10).Value = 240
$ws.Cells.Item(436, 11).Value = 4800
$ws.Cells.Item(436, 12).Value = 5000
$ws.Cells.Item(436, 13).Value = 4900
$ws.Cells.Item(436, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(436, 15).Value = "Región del Maule"
$ws.Cells.Item(436, 16).Value = 196
$ws.Cells.Item(436, 17).Value = 25
$ws.Cells.Item(436, 18).Value = "Hortaliza"
